$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.829.38'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.641.64'
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.97'
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("E8").Value = '  +0.59%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06433'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.40'
$ws.Range("E10").Value = '  +4.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07813'
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.294'
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.648.13'
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '1.871.29'
$ws.Range("E14").Value = '  +0.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5638'
$ws.Range("E15").Value = '  +2.67%  '
$ws.Range("D16").Value = '0.0₅7666'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.41'
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").Value = '25.863.79'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.22'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.369'
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.948'
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.116'
$ws.Range("E23").Value = '  +0.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.004'
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.804'
$ws.Range("E25").Value = '  -6.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.45'
$ws.Range("E26").Value = '  -1.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1248'
$ws.Range("E27").Value = '  +1.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.835'
$ws.Range("E28").Value = '  +0.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.53'
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.247'
$ws.Range("E30").Value = '  +0.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04943'
$ws.Range("E31").Value = '  +1.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.319'
$ws.Range("E32").Value = '  +2.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.252'
$ws.Range("E33").Value = '  +2.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.586'
$ws.Range("E34").Value = '  +3.14%  '
$ws.Range("E35").Value = '  +0.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9081'
$ws.Range("E36").Value = '  +1.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.576'
$ws.Range("E37").Value = '  +1.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5559'
$ws.Range("E38").Value = '  +0.86%  '
$ws.Range("D39").Value = '1.131.12'
$ws.Range("E39").Value = '  +1.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01568'
$ws.Range("E40").Value = '  +1.04%  '
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.543'
$ws.Range("E42").Value = '  -0.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8033'
$ws.Range("E43").Value = '  +0.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.52'
$ws.Range("E44").Value = '  +1.35%  '
$ws.Range("E45").Value = '  +0.88%  '
$ws.Range("E46").Value = '  -7.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.79'
$ws.Range("E47").Value = '  +2.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4276'
$ws.Range("E48").Value = '  -3.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.820'
$ws.Range("E49").Value = '  +3.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05036'
$ws.Range("E50").Value = '  -2.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9995'
$ws.Range("E51").Value = '  -0.24%  '
